# Have non-corner cases of choose_players working
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

# Widen column A to fit the new, longer label and drop the old best-fit autosize
$ws.Columns.Item(1).ColumnWidth = 25.17

# Row 3 / column A now reflects the locked state of the single-pick bot
$ws.Range("A3").Value = "singleBot(LOCKED)"

# New row for the 4th bot account
$ws.Range("A5").Value = "singleBot"
$ws.Range("B5").Value = "faiyamR004@gmail.com"
$ws.Range("C5").Value = "beatthestreak4"
$ws.Range("D5").Value = "beatthestreak4"
$ws.Range("E5").Value = "one selection"

# Hook up the mailto: hyperlink on the new email cell
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:faiyamR004@gmail.com")

# Re-apply the same Hyperlink look used by the other email cells (B2:B4)
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("E5").Select()
